$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 for the new "Burford Capital Limited" entry.
# This shifts the old row 3 (GLI Finance Limited) data down to row 4.
$ws.Rows("3:3").Insert()

# --- Row 2: update existing company record (now index "2") ---
$ws.Range("A2").Value = "Guernsey"
$ws.Range("B2").Value = "'2"
$ws.Range("C2").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("D2").Value = 0.1926
$ws.Range("E2").Value = 0.197
$ws.Range("G2").Value = 0.7924035952449986
$ws.Range("H2").Value = 0.7924035952449986
$ws.Range("I2").Value = 0.6648303856190201
$ws.Range("J2").Value = 0.5862498464436664
$ws.Range("K2").Value = 130.3
$ws.Range("L2").Value = 0.3777906639605683
$ws.Range("M2").Value = 9.119999999999999
$ws.Range("N2").Value = 0.004252144722118612
$ws.Range("O2").Value = 0.06999232540291635
$ws.Range("P2").Value = 9.119999999999999
$ws.Range("Q2").Value = 0.004252144722118612
$ws.Range("R2").Value = 0.06999232540291635
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 268.5
$ws.Range("V2").Value = 0.1251864975755315
$ws.Range("W2").Value = -0.06882966323535625
$ws.Range("X2").Value = 0.03503172946340932
$ws.Range("Y2").Value = -0.1038613926987656
$ws.Range("Z2").Value = 0.176956189362106
$ws.Range("AA2").Value = 0.04737166389179184
$ws.Range("AB2").Value = 0.02771141566137619
$ws.Range("AC2").Value = 0.01966024823041565
$ws.Range("AD2").Value = 723.9
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 723.9
$ws.Range("AG2").Value = 455.4
$ws.Range("AH2").Value = 0.2523442674382124
$ws.Range("AI2").Value = 0.2911204053727982
$ws.Range("AJ2").Value = 0.1751403738173986
$ws.Range("AK2").Value = 0.2053108516297732
$ws.Range("AL2").Value = 39.7
$ws.Range("AM2").Value = 39.7
$ws.Range("AN2").Value = 3.018765638031693
$ws.Range("AO2").Value = 5.775818639798489
$ws.Range("AP2").Value = 1.899082568807339
$ws.Range("AQ2").Value = 5.775818639798489

# --- Row 3: new company record (Burford Capital Limited) ---
$ws.Range("A3").Value = "Guernsey"
$ws.Range("B3").Value = "Burford Capital Limited (AIM:BUR)"
$ws.Range("C3").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("D3").Value = 0.31
$ws.Range("E3").Value = 0.197
$ws.Range("G3").Value = 0.8251811594202899
$ws.Range("H3").Value = 0.8251811594202899
$ws.Range("I3").Value = 0.6923309178743962
$ws.Range("J3").Value = 0.528668913275486
$ws.Range("K3").Value = 143.1
$ws.Range("L3").Value = 0.4320652173913043
$ws.Range("M3").Value = 9.119999999999999
$ws.Range("N3").Value = 0.004288939051918735
$ws.Range("O3").Value = 0.06373165618448637
$ws.Range("P3").Value = 9.119999999999999
$ws.Range("Q3").Value = 0.004288939051918735
$ws.Range("R3").Value = 0.06373165618448637
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 260.3
$ws.Range("V3").Value = 0.1224134687735139
$ws.Range("W3").Value = 0.09132099553286535
$ws.Range("X3").Value = 0.01939841129333628
$ws.Range("Y3").Value = 0.07192258423952907
$ws.Range("Z3").Value = 0.1792110816514258
$ws.Range("AA3").Value = 0.09474332778358369
$ws.Range("AB3").Value = 0.02496379999427911
$ws.Range("AC3").Value = 0.06977952778930457
$ws.Range("AD3").Value = 644.1
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 644.1
$ws.Range("AG3").Value = 383.8
$ws.Range("AH3").Value = 0.2324851109907959
$ws.Range("AI3").Value = 0.2723582392490169
$ws.Range("AJ3").Value = 0.1528961835710302
$ws.Range("AK3").Value = 0.1823624441699135
$ws.Range("AL3").Value = 39.7
$ws.Range("AM3").Value = 39.7
$ws.Range("AN3").Value = 2.685988323603003
$ws.Range("AO3").Value = 5.775818639798489
$ws.Range("AP3").Value = 1.600500417014179
$ws.Range("AQ3").Value = 5.775818639798489

# --- Row 4: GLI Finance Limited, shifted down, values refreshed ---
$ws.Range("A4").Value = "Guernsey"
$ws.Range("B4").Value = "GLI Finance Limited (AIM:GLIF)"
$ws.Range("C4").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("D4").Value = 0.07519999999999999
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = -12.8
$ws.Range("L4").Value = -0.9343065693430658
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 8.199999999999999
$ws.Range("V4").Value = 0.4456521739130435
$ws.Range("W4").Value = -0.2289803220035778
$ws.Range("X4").Value = 0.05066504763348235
$ws.Range("Y4").Value = -0.2796453696370602
$ws.Range("Z4").Value = 0.1356838664949985
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.03045903132847327
$ws.Range("AC4").Value = -0.03045903132847327
$ws.Range("AD4").Value = 79.8
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 79.8
$ws.Range("AG4").Value = 71.59999999999999
$ws.Range("AH4").Value = 0.8126272912423625
$ws.Range("AI4").Value = 0.6557107641741989
$ws.Range("AJ4").Value = 0.7955555555555555
$ws.Range("AK4").Value = 0.6308370044052863
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0

# Clear cells in row 4 that no longer have values after the refresh
# (previously populated for GLI Finance Limited at the old row 3 position)
$ws.Range("T4").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AO4").ClearContents()
$ws.Range("AP4").ClearContents()
$ws.Range("AQ4").ClearContents()
